$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.417.28"
$ws.Range("E2").Value = "  +4.00%  "

# Row 3
$ws.Range("D3").Value = "2.635.87"
$ws.Range("E3").Value = "  +2.20%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.65"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.50"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.44%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +4.16%  "

# Row 9
$ws.Range("D9").Value = "2.657.15"
$ws.Range("E9").Value = "  +2.95%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.83"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.37%  "

# Row 11
$ws.Range("E11").Value = "  +5.76%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.147"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +6.68%  "

# Row 13
$ws.Range("E13").Value = "  +4.19%  "

# Row 14
$ws.Range("D14").Value = "3.107.30"
$ws.Range("E14").Value = "  +2.38%  "

# Row 15
$ws.Range("D15").Value = "60.384.84"
$ws.Range("E15").Value = "  +4.05%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.12"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +6.96%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000138"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +5.45%  "

# Row 18
$ws.Range("D18").Value = "2.648.68"
$ws.Range("E18").Value = "  +1.72%  "

# Row 19
$ws.Range("E19").Value = "  +3.51%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "343.34"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.94%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.22%  "

# Row 22
$ws.Range("E22").Value = "  +3.54%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.08%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.21%  "

# Row 25
$ws.Range("E25").Value = "  +4.52%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.56%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.991"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.41"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.47%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0811"
$ws.Range("E29").Value = "  +11.94%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.08%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.71"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +4.98%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.18"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "159.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "19.19"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.32%  "

# Row 35
$ws.Range("E35").Value = "  +6.33%  "

# Row 36
$ws.Range("E36").Value = "  +8.80%  "

# Row 37
$ws.Range("E37").Value = "  +5.81%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.893"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +9.57%  "

# Row 39
$ws.Range("E39").Value = "  +7.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.47"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.60%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "299.66"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.93%  "

# Row 42
$ws.Range("E42").Value = "  +1.95%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.994"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0984"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +4.05%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.603"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.62%  "

# Row 46
$ws.Range("E46").Value = "  +2.78%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.44"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.32%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.18"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +16.03%  "

# Row 49
$ws.Range("E49").Value = "  +0.52%  "

# Row 50
$ws.Range("E50").Value = "  +4.20%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.71"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.99%  "
